$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 2.136171443821411
    "C2" = 0.1431366085257366
    "D2" = 0.09558933485272192
    "F2" = 2.235241546664398
    "G2" = 1.582972206320221
    "H2" = 1.403129638666357
    "J2" = 0.1806093174427268
    "L2" = 0.3960184363469921
    "N2" = 1.655404837749509
    "B3" = 2.028139294653442
    "C3" = 0.1270072910099884
    "D3" = 0.09539220275519611
    "F3" = 2.230732286596563
    "G3" = 1.573964234441632
    "H3" = 1.40532049914836
    "J3" = 0.1819304149090879
    "L3" = 0.3892893093184
    "N3" = 1.676187196226938
    "B4" = 1.962745785422783
    "C4" = 0.1170518062654651
    "D4" = 0.09528548125140546
    "F4" = 2.229240087674711
    "G4" = 1.569485028856064
    "H4" = 1.407389980972781
    "J4" = 0.1827985134864001
    "L4" = 0.3853269432126609
    "N4" = 1.689585237883247
    "B5" = 1.936334071356384
    "C5" = 0.1129817907936399
    "D5" = 0.0952456103371695
    "F5" = 2.228952577871112
    "G5" = 1.567923534066708
    "H5" = 1.408415285392579
    "J5" = 0.1831666005464037
    "L5" = 0.3837549213471618
    "N5" = 1.695205276515527
    "B6" = 1.931962739894402
    "C6" = 0.1123051783736457
    "D6" = 0.09523920892648619
    "F6" = 2.228924186946202
    "G6" = 1.567680164172074
    "H6" = 1.408596522760533
    "J6" = 0.1832285869636578
    "L6" = 0.3834964684697439
    "N6" = 1.696148154203254
    "B7" = 1.962388628687165
    "C7" = 0.1169969695412476
    "D7" = 0.09528492885996087
    "F7" = 2.229234912783213
    "G7" = 1.569462902664043
    "H7" = 1.407403071992135
    "J7" = 0.1828034196002708
    "L7" = 0.3853055694655154
    "N7" = 1.689660383070043
    "B8" = 2.098727478320313
    "C8" = 0.1375860148507115
    "D8" = 0.09551840280801116
    "F8" = 2.233421556912788
    "G8" = 1.579647574017713
    "H8" = 1.403734625851911
    "J8" = 0.1810530213018726
    "L8" = 0.3936631371924051
    "N8" = 1.662438204905987
    "B9" = 2.373523871562668
    "C9" = 0.1775500229222189
    "D9" = 0.09608910797196302
    "F9" = 2.251782433403108
    "G9" = 1.607997324705394
    "H9" = 1.402296189230469
    "J9" = 0.1780717398730509
    "L9" = 0.4113932729431724
    "N9" = 1.614118265443405
    "B10" = 2.579958862961689
    "C10" = 0.2066661731218176
    "D10" = 0.09657630261961003
    "F10" = 2.271496764255957
    "G10" = 1.63398337792475
    "H10" = 1.404762057382783
    "J10" = 0.1761557247661649
    "L10" = 0.4252353893227081
    "N10" = 1.581707905948335
    "B11" = 2.674861369040059
    "C11" = 0.2198598751995462
    "D11" = 0.09681250068269875
    "F11" = 2.281825211163778
    "G11" = 1.646936228270107
    "H11" = 1.40665200621666
    "J11" = 0.1753434911643943
    "L11" = 0.4317094522660909
    "N11" = 1.567635533358033
    "B12" = 2.710941234589313
    "C12" = 0.2248486416406763
    "D12" = 0.09690402146280519
    "F12" = 2.285932521020825
    "G12" = 1.652004691698124
    "H12" = 1.407478382304674
    "J12" = 0.175044446118882
    "L12" = 0.4341864388861438
    "N12" = 1.562403348625793
    "B13" = 2.703164468780358
    "C13" = 0.2237745509457056
    "D13" = 0.09688421869563157
    "F13" = 2.285039205775277
    "G13" = 1.650905823402923
    "H13" = 1.407295480673071
    "J13" = 0.1751084715517877
    "L13" = 0.433651847167809
    "N13" = 1.563525888589954
    "B14" = 2.677826831949346
    "C14" = 0.2202704524411843
    "D14" = 0.0968199886117489
    "F14" = 2.282159187728993
    "G14" = 1.647349932599781
    "H14" = 1.406717772907854
    "J14" = 0.1753187176320274
    "L14" = 0.4319127268769023
    "N14" = 1.567203138898088
    "B15" = 2.662325317374439
    "C15" = 0.2181231265605277
    "D15" = 0.09678091589108817
    "F15" = 2.28042065438413
    "G15" = 1.645193162250109
    "H15" = 1.406378332534928
    "J15" = 0.1754486100960122
    "L15" = 0.430850770785014
    "N15" = 1.569468162106643
    "B16" = 2.57377671357267
    "C16" = 0.2058028985491944
    "D16" = 0.09656115821567113
    "F16" = 2.270849194816876
    "G16" = 1.633159716200112
    "H16" = 1.404654022368874
    "J16" = 0.1762100002130076
    "L16" = 0.4248158529025403
    "N16" = 1.58264109226363
    "B17" = 2.519709205791742
    "C17" = 0.1982316677145661
    "D17" = 0.09643006335415194
    "F17" = 2.265326202967984
    "G17" = 1.626067977858668
    "H17" = 1.403793130833378
    "J17" = 0.1766922895154046
    "L17" = 0.4211589551231754
    "N17" = 1.590894318331086
    "B18" = 2.488704619512703
    "C18" = 0.1938720605287472
    "D18" = 0.09635603356276334
    "F18" = 2.262277529136583
    "G18" = 1.622095468882065
    "H18" = 1.403370262945089
    "J18" = 0.17697527829813
    "L18" = 0.4190722894535384
    "N18" = 1.595704526313594
    "B19" = 2.478223096709485
    "C19" = 0.1923951423905805
    "D19" = 0.09633120461516498
    "F19" = 2.26126726912814
    "G19" = 1.620768708837289
    "H19" = 1.403239496920634
    "J19" = 0.1770720535339798
    "L19" = 0.4183686485891087
    "N19" = 1.59734402398681
    "B20" = 2.525455100119473
    "C20" = 0.199038138711785
    "D20" = 0.09644387670467225
    "F20" = 2.265900882975657
    "G20" = 1.626811880609608
    "H20" = 1.403877290439993
    "J20" = 0.1766403706380828
    "L20" = 0.4215465118899573
    "N20" = 1.590009209067942
    "B21" = 2.685265251694148
    "C21" = 0.2212998921098688
    "D21" = 0.09683879829545106
    "F21" = 2.282999791003263
    "G21" = 1.648389940426483
    "H21" = 1.406884453677151
    "J21" = 0.1752567318118494
    "L21" = 0.4324228599886482
    "N21" = 1.566120414912903
    "B22" = 2.790539850379048
    "C22" = 0.2358060838079723
    "D22" = 0.0971090048304255
    "F22" = 2.295318491272994
    "G22" = 1.663445814519719
    "H22" = 1.409495122421589
    "J22" = 0.1744021610910167
    "L22" = 0.4396791920583638
    "N22" = 1.551071542582058
    "B23" = 2.734277123423567
    "C23" = 0.2280678087917636
    "D23" = 0.09696368888151596
    "F23" = 2.288638951099472
    "G23" = 1.655322722863417
    "H23" = 1.408042637653011
    "J23" = 0.1748537148022749
    "L23" = 0.4357928354476428
    "N23" = 1.55905175339926
    "B24" = 2.522857133477544
    "C24" = 0.1986735543734426
    "D24" = 0.09643762751888829
    "F24" = 2.265640675939053
    "G24" = 1.626475236336233
    "H24" = 1.403839017405062
    "J24" = 0.1766638253629402
    "L24" = 0.4213712484628616
    "N24" = 1.590409163516521
    "B25" = 2.298387212453122
    "C25" = 0.1667821747808489
    "D25" = 0.09592272894656517
    "F25" = 2.245724875653991
    "G25" = 1.599425908588898
    "H25" = 1.402067733234816
    "J25" = 0.1788300258675495
    "L25" = 0.4064534231831658
    "N25" = 1.626647709989514
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"